{"js": "// Update the first table's equation cells with the new values.\n// The replacement grid mirrors the table layout (20 rows x 5 columns),\n// read in row-major order, matching the target edit.\nconst newValues = [\n  [\n    \"23+30=\",\n    \"86-75=\",\n    \"15+72=\",\n    \"21+49=\",\n    \"90-49=\"\n  ],\n  [\n    \"43+41=\",\n    \"15+59=\",\n    \"4+67=\",\n    \"89-75=\",\n    \"32+37=\"\n  ],\n  [\n    \"86-34=\",\n    \"40+29=\",\n    \"16+28=\",\n    \"84-62=\",\n    \"81-57=\"\n  ],\n  [\n    \"58+25=\",\n    \"53-32=\",\n    \"29+33=\",\n    \"6+75=\",\n    \"7+26=\"\n  ],\n  [\n    \"23+70=\",\n    \"24+31=\",\n    \"20+15=\",\n    \"47-10=\",\n    \"84+12=\"\n  ],\n  [\n    \"91-71=\",\n    \"13+0=\",\n    \"26+45=\",\n    \"51-20=\",\n    \"62-35=\"\n  ],\n  [\n    \"67+20=\",\n    \"93-69=\",\n    \"25+41=\",\n    \"32+49=\",\n    \"20+58=\"\n  ],\n  [\n    \"49+14=\",\n    \"0+63=\",\n    \"36+29=\",\n    \"67+17=\",\n    \"17+53=\"\n  ],\n  [\n    \"84-26=\",\n    \"69-56=\",\n    \"61-49=\",\n    \"72-20=\",\n    \"46+15=\"\n  ],\n  [\n    \"79-76=\",\n    \"9+46=\",\n    \"41-28=\",\n    \"28+61=\",\n    \"47+51=\"\n  ],\n  [\n    \"10-6=\",\n    \"78-11=\",\n    \"0+84=\",\n    \"36+43=\",\n    \"88-1=\"\n  ],\n  [\n    \"60-18=\",\n    \"82-55=\",\n    \"86+2=\",\n    \"69-50=\",\n    \"67-1=\"\n  ],\n  [\n    \"43-3=\",\n    \"73+22=\",\n    \"96-56=\",\n    \"12-9=\",\n    \"82-47=\"\n  ],\n  [\n    \"30+25=\",\n    \"77-32=\",\n    \"19+13=\",\n    \"11+80=\",\n    \"49-4=\"\n  ],\n  [\n    \"73+11=\",\n    \"42+24=\",\n    \"91-21=\",\n    \"70+24=\",\n    \"57+29=\"\n  ],\n  [\n    \"10+86=\",\n    \"14+61=\",\n    \"24+52=\",\n    \"41-8=\",\n    \"63+12=\"\n  ],\n  [\n    \"69-45=\",\n    \"27+61=\",\n    \"49+38=\",\n    \"79-76=\",\n    \"95-52=\"\n  ],\n  [\n    \"13+44=\",\n    \"45-23=\",\n    \"48+17=\",\n    \"18+59=\",\n    \"56-34=\"\n  ],\n  [\n    \"17+25=\",\n    \"78-55=\",\n    \"0+60=\",\n    \"87-20=\",\n    \"91+0=\"\n  ],\n  [\n    \"83+11=\",\n    \"75+5=\",\n    \"49-38=\",\n    \"27+50=\",\n    \"19-15=\"\n  ]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nif (table.rowCount !== newValues.length) {\n  throw new Error(\n    `Unexpected row count: expected ${newValues.length}, got ${table.rowCount}`\n  );\n}\n\n// Rebuild each row preserving any existing columns beyond what we touch,\n// so only the equation text itself changes (formatting stays intact).\nconst updatedValues = table.values.map((row, r) =>\n  row.map((cell, c) => (newValues[r] && newValues[r][c] !== undefined ? newValues[r][c] : cell))\n);\n\ntable.values = updatedValues;\nawait context.sync();\n", "ps1": "# Update each cell of the first table with the new equation text.\n# Values are listed in row-major order (row 1 col 1..5, row 2 col 1..5, ...).\n$newValues = @(\n    \"23+30=\",\n    \"86-75=\",\n    \"15+72=\",\n    \"21+49=\",\n    \"90-49=\",\n    \"43+41=\",\n    \"15+59=\",\n    \"4+67=\",\n    \"89-75=\",\n    \"32+37=\",\n    \"86-34=\",\n    \"40+29=\",\n    \"16+28=\",\n    \"84-62=\",\n    \"81-57=\",\n    \"58+25=\",\n    \"53-32=\",\n    \"29+33=\",\n    \"6+75=\",\n    \"7+26=\",\n    \"23+70=\",\n    \"24+31=\",\n    \"20+15=\",\n    \"47-10=\",\n    \"84+12=\",\n    \"91-71=\",\n    \"13+0=\",\n    \"26+45=\",\n    \"51-20=\",\n    \"62-35=\",\n    \"67+20=\",\n    \"93-69=\",\n    \"25+41=\",\n    \"32+49=\",\n    \"20+58=\",\n    \"49+14=\",\n    \"0+63=\",\n    \"36+29=\",\n    \"67+17=\",\n    \"17+53=\",\n    \"84-26=\",\n    \"69-56=\",\n    \"61-49=\",\n    \"72-20=\",\n    \"46+15=\",\n    \"79-76=\",\n    \"9+46=\",\n    \"41-28=\",\n    \"28+61=\",\n    \"47+51=\",\n    \"10-6=\",\n    \"78-11=\",\n    \"0+84=\",\n    \"36+43=\",\n    \"88-1=\",\n    \"60-18=\",\n    \"82-55=\",\n    \"86+2=\",\n    \"69-50=\",\n    \"67-1=\",\n    \"43-3=\",\n    \"73+22=\",\n    \"96-56=\",\n    \"12-9=\",\n    \"82-47=\",\n    \"30+25=\",\n    \"77-32=\",\n    \"19+13=\",\n    \"11+80=\",\n    \"49-4=\",\n    \"73+11=\",\n    \"42+24=\",\n    \"91-21=\",\n    \"70+24=\",\n    \"57+29=\",\n    \"10+86=\",\n    \"14+61=\",\n    \"24+52=\",\n    \"41-8=\",\n    \"63+12=\",\n    \"69-45=\",\n    \"27+61=\",\n    \"49+38=\",\n    \"79-76=\",\n    \"95-52=\",\n    \"13+44=\",\n    \"45-23=\",\n    \"48+17=\",\n    \"18+59=\",\n    \"56-34=\",\n    \"17+25=\",\n    \"78-55=\",\n    \"0+60=\",\n    \"87-20=\",\n    \"91+0=\",\n    \"83+11=\",\n    \"75+5=\",\n    \"49-38=\",\n    \"27+50=\",\n    \"19-15=\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$expectedCount = $t.Rows.Count * $t.Columns.Count\nif ($newValues.Count -ne $expectedCount) {\n    throw \"Unexpected cell count: table has $expectedCount cells, but $($newValues.Count) replacement values were supplied.\"\n}\n\n$idx = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        $rng = $cell.Range\n        # Trim the trailing cell-mark characters (end-of-cell + paragraph mark)\n        # so only the visible text is replaced, preserving run formatting.\n        $rng.MoveEnd(1, -2) | Out-Null\n        $rng.Text = $newValues[$idx]\n        $idx++\n    }\n}\n"}
